# Commit: "version update edit product"
# - Update product image filename on row 4 (L4) from gio-qua-1.jpg to gio-qua-2.jpg
# - Move the active selection on the sheet to L6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the image name referenced by L4 (adds a new shared string entry,
# same as Excel would do when the cell text is edited in place).
$ws.Range("L4").Value = "gio-qua-2.jpg"

# Reflect the new selection left behind by the editing session.
$ws.Range("L6").Select() | Out-Null
